$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.175219416618347
$ws.Range("B1").Value = 2.400913238525391
$ws.Range("D1").Value = 2.352333068847656
$ws.Range("E1").Value = 1.20781409740448
